# Add data from 2020 Aug
# ------------------------------------------------------------------
# 1. Append 6 new rows (10-15) of observations to delta_cep_2020,
#    preserving number formats by copying the last existing data row
#    (row 9) down, then overwriting with the real values/formulas.
# 2. Extend the scatter-chart series ("2020") so it covers the new
#    rows as well.
# 3. Update sheet-view/window state (active sheet moved from mu_cep
#    to delta_cep, selections moved, window height changed).
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. New rows on delta_cep_2020
# ---------------------------------------------------------------
$ws2020 = $wb.Worksheets.Item("delta_cep_2020")

# Clone the formatting/formulas of row 9 five more times so rows
# 10..15 inherit the same styles (date format on A, number format on F)
# and relative-formula shape.
for ($i = 10; $i -le 15; $i++) {
    $ws2020.Rows("9:9").Copy()
    $ws2020.Rows("$($i):$($i)").Insert(-4121, 0)
}

# Row 10
$ws2020.Range("A10").Value = 43881
$ws2020.Range("B10").Formula = "=A10-delta_cep!A10+delta_cep!B10"
$ws2020.Range("C10").Value = 18
$ws2020.Range("D10").Value = 46
$ws2020.Range("E10").Formula = "=(B10*1440+C10*60+D10)/1440"
$ws2020.Range("F10").Value = 3.7
$ws2020.Range("G10").Formula = "=E10-5.366*H10"
$ws2020.Range("H10").Value = 159

# Row 11 (minute value kept as text "02")
$ws2020.Range("A11").Value = 44055
$ws2020.Range("B11").Formula = "=A11-delta_cep!A12+delta_cep!B12"
$ws2020.Range("C11").Value = 22
$ws2020.Range("D11").NumberFormat = "@"
$ws2020.Range("D11").Value = "02"
$ws2020.Range("E11").Formula = "=(B11*1440+C11*60+D11)/1440"
$ws2020.Range("F11").Value = 4.2
$ws2020.Range("G11").Formula = "=E11-5.366*H11"
$ws2020.Range("H11").Value = 191

# Row 12
$ws2020.Range("A12").Value = 44056
$ws2020.Range("B12").Formula = "=A12-delta_cep!A13+delta_cep!B13"
$ws2020.Range("C12").Value = 21
$ws2020.Range("D12").Value = 28
$ws2020.Range("E12").Formula = "=(B12*1440+C12*60+D12)/1440"
$ws2020.Range("F12").Value = 3.5
$ws2020.Range("G12").Formula = "=E12-5.366*H12"
$ws2020.Range("H12").Value = 191

# Row 13
$ws2020.Range("A13").Value = 44057
$ws2020.Range("B13").Formula = "=A13-delta_cep!A14+delta_cep!B14"
$ws2020.Range("C13").Value = 21
$ws2020.Range("D13").Value = 31
$ws2020.Range("E13").Formula = "=(B13*1440+C13*60+D13)/1440"
$ws2020.Range("F13").Value = 3.6
$ws2020.Range("G13").Formula = "=E13-5.366*H13"
$ws2020.Range("H13").Value = 191

# Row 14
$ws2020.Range("A14").Value = 44058
$ws2020.Range("B14").Formula = "=A14-delta_cep!A15+delta_cep!B15"
$ws2020.Range("C14").Value = 23
$ws2020.Range("D14").Value = 24
$ws2020.Range("E14").Formula = "=(B14*1440+C14*60+D14)/1440"
$ws2020.Range("F14").Value = 3.8
$ws2020.Range("G14").Formula = "=E14-5.366*H14"
$ws2020.Range("H14").Value = 192

# Row 15
$ws2020.Range("A15").Value = 44059
$ws2020.Range("B15").Formula = "=A15-delta_cep!A16+delta_cep!B16"
$ws2020.Range("C15").Value = 20
$ws2020.Range("D15").Value = 47
$ws2020.Range("E15").Formula = "=(B15*1440+C15*60+D15)/1440"
$ws2020.Range("F15").Value = 4.2
$ws2020.Range("G15").Formula = "=E15-5.366*H15"
$ws2020.Range("H15").Value = 192

# Selection left where the user clicked last on this sheet
$ws2020.Range("G21").Select()

# ---------------------------------------------------------------
# 2. Extend the "2020" series on the chart living on delta_cep
# ---------------------------------------------------------------
$wsDelta = $wb.Worksheets.Item("delta_cep")
$chart = $wsDelta.ChartObjects(1).Chart
$series2020 = $chart.SeriesCollection(4)
$series2020.XValues = "=delta_cep_2020!`$G`$2:`$G`$15"
$series2020.Values = "=delta_cep_2020!`$F`$2:`$F`$15"

# ---------------------------------------------------------------
# 3. Sheet / window state: active sheet is now delta_cep (was
#    mu_cep), with a fresh selection near the chart; window height
#    shrinks slightly.
# ---------------------------------------------------------------
$wsMu = $wb.Worksheets.Item("mu_cep")
$wsMu.Range("F29").Select()

$wsDelta.Activate()
$wsDelta.Range("W14").Select()

$excel.ActiveWindow.Height = 9885
